# Update the dSF column (F) values for the musgrove_joe dataset.
# This reflects a "repull" of the underlying data (mean calc) that changed
# the dSF figures for most rows (row 1 is the header).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = -1
    3  = -2
    4  = 4
    7  = 1
    8  = 2
    10 = -4
    11 = 3
    12 = 2
    14 = -3
    15 = 2
    16 = -3
    17 = -1
    18 = -2
    19 = 3
    20 = 1
    21 = -2
    22 = -5
    23 = 3
    24 = 5
    25 = 6
    27 = 2
    28 = -2
    29 = -2
    30 = 2
    31 = 4
    32 = 3
    33 = 7
    35 = -2
    36 = -2
    37 = -2
}

foreach ($row in $newValues.Keys) {
    $ws.Range("F$row").Value = $newValues[$row]
}
